# Append three new entries ("Réalisation") to the work-journal table
# (Tableau1 on the active sheet), matching the data the table grows to
# (A1:F89 -> A1:F92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$newEntries = @(
    @{ Row = 90; Date = 44705; Type = "Réalisation"; Duree = 0.5;  Desc = "Revue de l'ortographe sur le site et la doc" },
    @{ Row = 91; Date = 44705; Type = "Réalisation"; Duree = 0.75; Desc = "Recherche d'articles dans une seule catégorie" },
    @{ Row = 92; Date = 44705; Type = "Réalisation"; Duree = 0.25; Desc = "Documentation du travail de cet après-midi" }
)

foreach ($entry in $newEntries) {
    $lo.ListRows.Add() | Out-Null

    $dateCell = $ws.Cells.Item($entry.Row, 1)
    # Reuse the date-formatted style already used by the row above instead of
    # minting a brand-new number format, so the new cells render as dates.
    $ws.Range($ws.Cells.Item($entry.Row - 1, 1), $ws.Cells.Item($entry.Row - 1, 1)).Copy() | Out-Null
    $dateCell.PasteSpecial(-4122) | Out-Null

    $dateCell.Value = $entry.Date
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Type
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Duree
    $ws.Cells.Item($entry.Row, 4).Value = $entry.Desc
}

$excel.CutCopyMode = 0

$ws.Range("D92").Select() | Out-Null
